$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Fix existing shared-formula group so F22 becomes its own formula
# (row 20/21 keep the shared formula, row 22 is explicit)
# ---------------------------------------------------------------
$ws.Range("F22").Formula = "=(E22-E23)*(1-F`$18)"

# ---------------------------------------------------------------
# New block of data below the existing tables: rows 28-34
# Mirrors the "analytical" block at rows 17-23 but using only two
# filament diameters (0.3 and 0.2) for a 3d / small-length case.
# ---------------------------------------------------------------

# Row 28 - header row
$ws.Range("C28").Value = "Diameter"
$ws.Range("D28").Value = "Volume"
$ws.Range("E28").Value = "analytical"
$ws.Range("F28").Value = "Vf"
$ws.Range("G28").Value = "# particles"
$ws.Range("H28").Value = "Probability"

# Row 29 - GC (base) row
$ws.Range("B29").Value = "GC"
$ws.Range("C29").Value = 1
$ws.Range("C29").Style = "Input"
$ws.Range("D29").Formula = "=C29^2*PI()/4"
$ws.Range("F29").Value = 0.7
$ws.Range("F29").Style = "Input"
$ws.Range("G29").Formula = "=(F29/D29)/(F`$3/D`$3)"
$ws.Range("H29").Formula = "=G29/`$G`$34"
$ws.Range("H29").NumberFormat = "0.00000"

# Rows 30-31 - fil_1 / fil_2
$ws.Range("B30").Value = "fil_1"
$ws.Range("B31").Value = "fil_2"
$ws.Range("C30").Value = 0.3
$ws.Range("C31").Value = 0.2
$ws.Range("C30:C31").Style = "Input"

$ws.Range("D30:D31").Formula = "=C30^2*PI()/4"

$ws.Range("E30").Formula = "=(C30/`$C`$30)^`$E`$1"
$ws.Range("E31").Formula = "=(C31/`$C`$30)^`$E`$1"

$ws.Range("F30").Formula = "=(E30-E31)*(1-F`$29)"
$ws.Range("F31").Formula = "=(E31-E32)*(1-F`$29)"

$ws.Range("G30").Formula = "=(F30/D30)/(F`$3/D`$3)"
$ws.Range("G31").Formula = "=(F31/D31)/(F`$3/D`$3)"

$ws.Range("H30:H31").Formula = "=G30/`$G`$34"
$ws.Range("H30:H31").NumberFormat = "0.00000"

# Rows 32-33 - blank spacer rows (styled like the column above/right)
$ws.Range("C32").Style = "Input"
$ws.Range("H32").NumberFormat = "0.00000"
$ws.Range("C33").Style = "Input"
$ws.Range("H33").NumberFormat = "0.00000"

# Row 34 - totals
$ws.Range("G34").Formula = "=SUM(G29:G33)"
$ws.Range("G34").Font.Name = "Calibri"
$ws.Range("H34").Formula = "=SUM(H29:H33)"
$ws.Range("H34").NumberFormat = "0.00000"

# ---------------------------------------------------------------
# Update the active selection to the last cell edited
# ---------------------------------------------------------------
$ws.Range("I34").Select()
